# fancy_table.docx edit:
#  - widen gridCol #4 and #7 (1132 -> 1181 twips)
#  - bump the column-header row's auto height (614 -> 615 twips)
#  - bold the column-header labels in that row (leaving the superscript
#    footnote-marker runs, e.g. the "1" after "OR", un-bolded)

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- widen the 4th and 7th table-grid columns -----------------------------
# Word's Column.Width is expressed in points; 1 twip = 1/20 point.
$t.Columns.Item(4).Width = 1181 / 20
$t.Columns.Item(7).Width = 1181 / 20

# --- bump the header row (row 2 - "Characteristic"/"OR"/...) height -------
# (row already uses wdRowHeightAuto; just nudge the stored twip value)
$headerRow = $t.Rows.Item(2)
$headerRow.Height = 615 / 20

# --- bold the header-row labels, skipping the superscript "1" markers -----
$labels = @("Characteristic", "OR", "95% CI", "p-value", "HR", "95% CI", "p-value")
for ($c = 1; $c -le $t.Columns.Count; $c++) {
    $cell = $t.Cell(2, $c)
    $cellRange = $cell.Range
    $label = $labels[$c - 1]
    $labelRange = $d.Range($cellRange.Start, $cellRange.Start + $label.Length)
    $labelRange.Bold = 1
}
